# Updated cryptos list on Thu Nov 30 15:58:13 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force cells whose content looks like a plain number (single decimal point)
    # to stay as literal text, matching the original inline-string content.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.590.93"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.029.94"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "226.25"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.29%  "

# Row 7 - Solana
Set-TextValue "D7" "59.55"
$ws.Range("E7").Value = "  -1.28%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.62%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0824"
$ws.Range("E10").Value = "  +1.95%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.32%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.329.22"
$ws.Range("E12").Value = "  -0.16%  "

# Row 13 - Chainlink
Set-TextValue "D13" "14.37"
$ws.Range("E13").Value = "  -1.71%  "

# Row 14 - Avalanche
Set-TextValue "D14" "21.01"
$ws.Range("E14").Value = "  -0.96%  "

# Row 15 - Polkadot
Set-TextValue "D15" "5.50"
$ws.Range("E15").Value = "  +4.88%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.769"
$ws.Range("E16").Value = "  +1.72%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.033.62"
$ws.Range("E17").Value = "  -0.27%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.542.11"
$ws.Range("E18").Value = "  -0.85%  "

# Row 20 - Litecoin
Set-TextValue "D20" "69.23"
$ws.Range("E20").Value = "  -0.63%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("E21").Value = "  -0.59%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "223.64"
$ws.Range("E22").Value = "  -0.39%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.25%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.26"
$ws.Range("E25").Value = "  +2.99%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.37"
$ws.Range("E26").Value = "  +2.13%  "

# Row 27 - Monero
Set-TextValue "D27" "167.53"
$ws.Range("E27").Value = "  +1.33%  "

# Row 28 - Kaspa
Set-TextValue "D28" "0.127"
$ws.Range("E28").Value = "  -1.33%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.92%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -1.91%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +0.35%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  +8.90%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -1.88%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0605"
$ws.Range("E34").Value = "  +0.19%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("E35").Value = "  -1.62%  "

# Row 36 - THORChain
$ws.Range("E36").Value = "  +2.50%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +3.33%  "

# Row 38 - RenderToken
Set-TextValue "D38" "3.39"
$ws.Range("E38").Value = "  +4.67%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  -0.03%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "17.88"
$ws.Range("E40").Value = "  +7.86%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.521.30"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42 - Aave
Set-TextValue "D42" "96.58"
$ws.Range("E42").Value = "  -0.40%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -1.33%  "

# Row 44 - HuobiToken
$ws.Range("E44").Value = "  +2.45%  "

# Row 45 - was Cronos, now FTXToken (rows 45/46 swapped)
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D45" "4.22"
$ws.Range("E45").Value = "  +5.63%  "

# Row 46 - was FTXToken, now Cronos
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0905"
$ws.Range("E46").Value = "  -1.49%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  -0.56%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  -0.09%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  -1.28%  "

# Row 50 - FraxShare
Set-TextValue "D50" "7.05"
$ws.Range("E50").Value = "  -0.03%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.220.62"
$ws.Range("E51").Value = "  -0.12%  "
